# Auto-generated edit script: apply Ultros_Profits market-price updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 824.7273
$ws.Range("I19").Value = 1170.25
$ws.Range("K19").Value = 1170.25
$ws.Range("M19").Value = -995.25

$ws.Range("H80").Value = 1651.7428
$ws.Range("I80").Value = 714.05554
$ws.Range("J80").Value = 2644.5881
$ws.Range("K80").Value = 2142.16662
$ws.Range("L80").Value = 7933.7643
$ws.Range("M80").Value = -1144.16662
$ws.Range("N80").Value = -9929.764299999999

$ws.Range("H83").Value = 1651.7428
$ws.Range("I83").Value = 714.05554
$ws.Range("J83").Value = 2644.5881
$ws.Range("K83").Value = 6426.49986
$ws.Range("L83").Value = 23801.2929
$ws.Range("M83").Value = -1434.49986
$ws.Range("N83").Value = -33785.2929

$ws.Range("H100").Value = 6259.1816
$ws.Range("I100").Value = 3968
$ws.Range("J100").Value = 8168.5
$ws.Range("K100").Value = 3968
$ws.Range("L100").Value = 8168.5
$ws.Range("M100").Value = -3427
$ws.Range("N100").Value = -9250.5

$ws.Range("H132").Value = 1072.1428
$ws.Range("I132").Value = 924.4722
$ws.Range("K132").Value = 2773.4166
$ws.Range("M132").Value = -243.4166

$ws.Range("H137").Value = 3717.516
$ws.Range("I137").Value = 2410.0952
$ws.Range("J137").Value = 6463.1
$ws.Range("K137").Value = 7230.285600000001
$ws.Range("L137").Value = 19389.3
$ws.Range("M137").Value = -4680.285600000001
$ws.Range("N137").Value = -24489.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 3175
$ws.Range("I25").Value = 3175
$ws.Range("K25").Value = 3175
$ws.Range("M25").Value = -2773

$ws.Range("H32").Value = 6978.58
$ws.Range("I32").Value = 5772.1895
$ws.Range("K32").Value = 5772.1895
$ws.Range("M32").Value = -5485.1895

$ws.Range("H45").Value = 3450.7273
$ws.Range("I45").Value = 1216
$ws.Range("K45").Value = 1216
$ws.Range("M45").Value = -839

$ws.Range("H98").Value = 58632.668
$ws.Range("J98").Value = 58632.668
$ws.Range("L98").Value = 58632.668
$ws.Range("N98").Value = -64622.668

$ws.Range("H102").Value = 2262.7058
$ws.Range("I102").Value = 1979.2667
$ws.Range("K102").Value = 1979.2667
$ws.Range("M102").Value = -357.2666999999999

$ws.Range("H104").Value = 39999
$ws.Range("J104").Value = 39999
$ws.Range("L104").Value = 39999
$ws.Range("N104").Value = -46987

$ws.Range("H107").Value = 48999.5
$ws.Range("J107").Value = 48999.5
$ws.Range("L107").Value = 48999.5
$ws.Range("N107").Value = -56679.5

$ws.Range("H122").Value = 3896.2917
$ws.Range("I122").Value = 2969.4375
$ws.Range("K122").Value = 8908.3125
$ws.Range("M122").Value = -6458.3125

$ws.Range("H132").Value = 4442.7905
$ws.Range("I132").Value = 3975.3845
$ws.Range("K132").Value = 11926.1535
$ws.Range("M132").Value = -9396.1535

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2734.625
$ws.Range("I86").Value = 3126.111
$ws.Range("J86").Value = 2231.2856
$ws.Range("K86").Value = 3126.111
$ws.Range("L86").Value = 2231.2856
$ws.Range("M86").Value = -2003.111
$ws.Range("N86").Value = -4477.2856

$ws.Range("H89").Value = 2734.625
$ws.Range("I89").Value = 3126.111
$ws.Range("J89").Value = 2231.2856
$ws.Range("K89").Value = 15630.555
$ws.Range("L89").Value = 11156.428
$ws.Range("M89").Value = -10014.555
$ws.Range("N89").Value = -22388.428

$ws.Range("H92").Value = 58000
$ws.Range("J92").Value = 58000
$ws.Range("L92").Value = 58000
$ws.Range("N92").Value = -62992

$ws.Range("H99").Value = 1999.6666
$ws.Range("I99").Value = 1999.6666
$ws.Range("K99").Value = 1999.6666
$ws.Range("M99").Value = -501.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3530.1853
$ws.Range("I31").Value = 2684.7812
$ws.Range("J31").Value = 4759.864
$ws.Range("K31").Value = 2684.7812
$ws.Range("L31").Value = 4759.864
$ws.Range("M31").Value = -2389.7812
$ws.Range("N31").Value = -5349.864

$ws.Range("H34").Value = 3530.1853
$ws.Range("I34").Value = 2684.7812
$ws.Range("J34").Value = 4759.864
$ws.Range("K34").Value = 2684.7812
$ws.Range("L34").Value = 4759.864
$ws.Range("M34").Value = -2482.7812
$ws.Range("N34").Value = -5163.864

$ws.Range("H62").Value = 6939.8
$ws.Range("I62").Value = 6424.75
$ws.Range("J62").Value = 9000
$ws.Range("K62").Value = 6424.75
$ws.Range("L62").Value = 9000
$ws.Range("M62").Value = -5800.75
$ws.Range("N62").Value = -10248

$ws.Range("H65").Value = 6939.8
$ws.Range("I65").Value = 6424.75
$ws.Range("J65").Value = 9000
$ws.Range("K65").Value = 32123.75
$ws.Range("L65").Value = 45000
$ws.Range("M65").Value = -29003.75
$ws.Range("N65").Value = -51240

$ws.Range("H107").Value = 6777.8823
$ws.Range("I107").Value = 857
$ws.Range("K107").Value = 857
$ws.Range("M107").Value = 1063

$ws.Range("H132").Value = 3796.5789
$ws.Range("I132").Value = 2931.1333
$ws.Range("K132").Value = 8793.3999
$ws.Range("M132").Value = -6263.3999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1269.6
$ws.Range("I8").Value = 1269.6
$ws.Range("K8").Value = 3808.8
$ws.Range("M8").Value = -3669.8

$ws.Range("H22").Value = 47620210
$ws.Range("I22").Value = 832.6667
$ws.Range("K22").Value = 2498.0001
$ws.Range("M22").Value = -2329.0001

$ws.Range("H27").Value = 47620210
$ws.Range("I27").Value = 832.6667
$ws.Range("K27").Value = 2498.0001
$ws.Range("M27").Value = -2396.0001

$ws.Range("H29").Value = 22222908
$ws.Range("I29").Value = 1134.125
$ws.Range("J29").Value = 47619220
$ws.Range("K29").Value = 3402.375
$ws.Range("L29").Value = 142857660
$ws.Range("M29").Value = -3125.375
$ws.Range("N29").Value = -142858214

$ws.Range("H122").Value = 2388.25
$ws.Range("J122").Value = 2581.875
$ws.Range("L122").Value = 23236.875
$ws.Range("N122").Value = -28136.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1600
$ws.Range("I97").Value = 1457.1428
$ws.Range("J97").Value = 1850
$ws.Range("K97").Value = 1457.1428
$ws.Range("L97").Value = 1850
$ws.Range("M97").Value = -961.1428000000001
$ws.Range("N97").Value = -2842

$ws.Range("H132").Value = 4631.452
$ws.Range("I132").Value = 4842.3057
$ws.Range("J132").Value = 3366.3333
$ws.Range("K132").Value = 14526.9171
$ws.Range("L132").Value = 10098.9999
$ws.Range("M132").Value = -11996.9171
$ws.Range("N132").Value = -15158.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8350.718999999999
$ws.Range("I46").Value = 2413.6667
$ws.Range("J46").Value = 8964.896000000001
$ws.Range("K46").Value = 2413.6667
$ws.Range("L46").Value = 8964.896000000001
$ws.Range("M46").Value = -2225.6667
$ws.Range("N46").Value = -9340.896000000001

$ws.Range("H48").Value = 19498
$ws.Range("J48").Value = 24247
$ws.Range("L48").Value = 24247
$ws.Range("N48").Value = -25569

$ws.Range("H55").Value = 2188.0833
$ws.Range("I55").Value = 2425.8
$ws.Range("K55").Value = 2425.8
$ws.Range("M55").Value = -2252.8

$ws.Range("H97").Value = 24499
$ws.Range("J97").Value = 24499
$ws.Range("L97").Value = 24499
$ws.Range("N97").Value = -26481

$ws.Range("H100").Value = 77385.734
$ws.Range("I100").Value = 124731.78
$ws.Range("K100").Value = 124731.78
$ws.Range("M100").Value = -124190.78

$ws.Range("H122").Value = 6354.615
$ws.Range("I122").Value = 3826.25
$ws.Range("K122").Value = 11478.75
$ws.Range("M122").Value = -9028.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 1316.6666
$ws.Range("J7").Value = 950
$ws.Range("L7").Value = 950
$ws.Range("N7").Value = -1176

$ws.Range("H52").Value = 3768566
$ws.Range("I52").Value = 6012806
$ws.Range("J52").Value = 28166.334
$ws.Range("K52").Value = 6012806
$ws.Range("L52").Value = 28166.334
$ws.Range("M52").Value = -6012580
$ws.Range("N52").Value = -28618.334

$ws.Range("H122").Value = 6983.9165
$ws.Range("I122").Value = 6225.875
$ws.Range("K122").Value = 18677.625
$ws.Range("M122").Value = -16227.625

$ws.Range("H126").Value = 2290.875
$ws.Range("I126").Value = 1983.4286
$ws.Range("K126").Value = 5950.2858
$ws.Range("M126").Value = -3480.2858

$ws.Range("H136").Value = 3595.9302
$ws.Range("I136").Value = 3204.8462
$ws.Range("J136").Value = 4194.0586
$ws.Range("K136").Value = 9614.5386
$ws.Range("L136").Value = 12582.1758
$ws.Range("M136").Value = -7064.5386
$ws.Range("N136").Value = -17682.1758

Write-Host "Applied 220 cell updates"